$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.595.40'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '2.433.36'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.03'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.98'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').Value = '2.445.86'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  -4.52%  '
$ws.Range('E12').Value = '  -4.51%  '
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('D14').Value = '2.864.80'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').Value = '57.490.02'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('E17').Value = '  -3.16%  '
$ws.Range('D18').Value = '2.439.07'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.44'
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.10'
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.76'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.22'
$ws.Range('E24').Value = '  -1.89%  '
$ws.Range('E25').Value = '  -1.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.23'
$ws.Range('E28').Value = '  -2.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.75'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0720'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.22'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.65'
$ws.Range('E32').Value = '  -3.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.65'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.90'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.21'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  -2.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.760'
$ws.Range('E41').Value = '  -4.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '271.27'
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.85'
$ws.Range('E44').Value = '  -2.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.578'
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.02'
$ws.Range('E47').Value = '  -5.25%  '
$ws.Range('E48').Value = '  -2.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.09'
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0208'
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.53'
$ws.Range('E51').Value = '  -3.23%  '
